$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy header style from an existing header cell (e.g. H1) to I1 and J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF)
$values = @{
    2  = @(6, 7)
    3  = @(7, 8)
    4  = @(4, 4)
    5  = @(8, 8)
    6  = @(8, 8)
    7  = @(5, 5)
    8  = @(6, 6)
    9  = @(7, 7)
    10 = @(7, 8)
    11 = @(12, 12)
    12 = @(6, 6)
    13 = @(5, 5)
    14 = @(8, 8)
    15 = @(5, 6)
    16 = @(10, 10)
    17 = @(5, 5)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(5, 6)
    21 = @(6, 6)
    22 = @(7, 7)
    23 = @(6, 6)
    24 = @(7, 7)
    25 = @(7, 7)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
